$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("singleInputData")
$ws2 = $wb.Worksheets.Item("doubleInputData")

# --- singleInputData (sheet1): insert two new rows before the old row 6 ---
$ws1.Rows("6:7").Insert()
$ws1.Range("A6").Value = "15fg"
$ws1.Range("B6").Value = "15fg"
$ws1.Range("A7").Value = "M1ATG"
$ws1.Range("B7").Value = "M1ATG"

# --- doubleInputData (sheet2): insert a new row before the old row 6 ---
$ws2.Rows("6:6").Insert()
$ws2.Range("A6").Value = 12
$ws2.Range("B6").Value = 73
$ws2.Range("C6").Value = 85

# --- selection / active sheet bookkeeping ---
$ws2.Range("E13").Select()
$ws1.Activate()
